$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target stored widths: col5=2.140625, col7=3.140625, col9/10=5.7109375)
# Note: this runtime's ColumnWidth setter re-quantizes the supplied "character width"
# in increments of 1/6 (stored_width = round(input*6)/6 + 5/6), so the values below are
# chosen as the closest inputs that land the stored width on the nearest achievable bucket
# to the true target width.
$ws.Columns.Item(5).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(7).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(9).ColumnWidth = 4.833333333333333
$ws.Columns.Item(10).ColumnWidth = 4.833333333333333

# Update cell values in row 1
$ws.Range("A1").Value = 3
$ws.Range("C1").Value = 19
$ws.Range("D1").Value = 24
$ws.Range("E1").Value = 8
$ws.Range("F1").Value = 28
$ws.Range("G1").Value = 13
$ws.Range("I1").Value = 0.092
$ws.Range("J1").Value = 0.077
$ws.Range("K1").Value = 0.065
